# [FEAT] Ordering columns and check empties
#
# Column C ("Coluna C") is moved before column B ("Coluna B"), rows 2/3
# are re-ordered, and the (now) empty/duplicate cell is cleared while a
# new "Linha BB" value is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: swap "Coluna B" and "Coluna C" ---
$ws.Range("B1").Value = "Coluna C"
$ws.Range("C1").Value = "Coluna B"

# --- Row 2 (was row 3's A/B data) ---
$ws.Range("A2").Value = "Linha A3"
$ws.Range("B2").Value = "Linha C3"
$ws.Range("C2").Value = ""

# --- Row 3 (was row 2's A/B data), plus new value in C3 ---
$ws.Range("A3").Value = "Linha A2"
$ws.Range("B3").Value = "Linha C2"
$ws.Range("C3").Value = "Linha BB"

# --- Row 4: keep A4, swap B4/C4 content (Coluna C/B swap) ---
$ws.Range("A4").Value = "Linha A4"
$ws.Range("B4").Value = "Linha C4"
$ws.Range("C4").Value = "Linha B4"

# --- View: select C1:C4 and zoom in ---
$ws.Select()
$ws.Range("C1:C4").Select()
$excel.ActiveWindow.Zoom = 265
